# Auto-generated edit script: apply numeric corrections to Ramuh_Profits workbook
# as described by the commit diff (profit/price recompute across several sheets).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132: H132=3283.215, I132=3374.9077, J132=3070.3572, K132=10124.7231, L132=9211.071599999999, M132=-7594.723100000001, N132=-14271.0716
$ws.Range("H132").Value = 3283.215
$ws.Range("I132").Value = 3374.9077
$ws.Range("J132").Value = 3070.3572
$ws.Range("K132").Value = 10124.7231
$ws.Range("L132").Value = 9211.071599999999
$ws.Range("M132").Value = -7594.723100000001
$ws.Range("N132").Value = -14271.0716

$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32=5099.7407, I32=4267.7827, K32=4267.7827, M32=-3980.7827
$ws.Range("H32").Value = 5099.7407
$ws.Range("I32").Value = 4267.7827
$ws.Range("K32").Value = 4267.7827
$ws.Range("M32").Value = -3980.7827
# Row 63: H63=1790.4615, J63=1870.75, L63=1870.75, N63=-3242.75
$ws.Range("H63").Value = 1790.4615
$ws.Range("J63").Value = 1870.75
$ws.Range("L63").Value = 1870.75
$ws.Range("N63").Value = -3242.75
# Row 66: H66=1790.4615, J66=1870.75, L66=9353.75, N66=-16217.75
$ws.Range("H66").Value = 1790.4615
$ws.Range("J66").Value = 1870.75
$ws.Range("L66").Value = 9353.75
$ws.Range("N66").Value = -16217.75
# Row 74: H74=3141.7673, I74=772.7857, J74=7563.8667, K74=772.7857, L74=7563.8667, M74=101.2143, N74=-9311.866699999999
$ws.Range("H74").Value = 3141.7673
$ws.Range("I74").Value = 772.7857
$ws.Range("J74").Value = 7563.8667
$ws.Range("K74").Value = 772.7857
$ws.Range("L74").Value = 7563.8667
$ws.Range("M74").Value = 101.2143
$ws.Range("N74").Value = -9311.866699999999
# Row 77: H77=3141.7673, I77=772.7857, J77=7563.8667, K77=3863.9285, L77=37819.3335, M77=504.0715, N77=-46555.3335
$ws.Range("H77").Value = 3141.7673
$ws.Range("I77").Value = 772.7857
$ws.Range("J77").Value = 7563.8667
$ws.Range("K77").Value = 3863.9285
$ws.Range("L77").Value = 37819.3335
$ws.Range("M77").Value = 504.0715
$ws.Range("N77").Value = -46555.3335
# Row 97: H97=614.6429000000001, I97=614.6429000000001, K97=614.6429000000001, M97=-118.6429000000001
$ws.Range("H97").Value = 614.6429000000001
$ws.Range("I97").Value = 614.6429000000001
$ws.Range("K97").Value = 614.6429000000001
$ws.Range("M97").Value = -118.6429000000001
# Row 102: H102=1129.909, I102=1042.9, K102=1042.9, M102=579.0999999999999
$ws.Range("H102").Value = 1129.909
$ws.Range("I102").Value = 1042.9
$ws.Range("K102").Value = 1042.9
$ws.Range("M102").Value = 579.0999999999999
# Row 132: H132=6021.976, I132=4030.8857, J132=15977.429, K132=12092.6571, L132=47932.287, M132=-9562.6571, N132=-52992.287
$ws.Range("H132").Value = 6021.976
$ws.Range("I132").Value = 4030.8857
$ws.Range("J132").Value = 15977.429
$ws.Range("K132").Value = 12092.6571
$ws.Range("L132").Value = 47932.287
$ws.Range("M132").Value = -9562.6571
$ws.Range("N132").Value = -52992.287

$ws = $wb.Worksheets.Item("BSM")
# Row 94: H94=17858138, I94=33334034, J94=1336.5385, K94=33334034, L94=1336.5385, M94=-33333583, N94=-2238.5385
$ws.Range("H94").Value = 17858138
$ws.Range("I94").Value = 33334034
$ws.Range("J94").Value = 1336.5385
$ws.Range("K94").Value = 33334034
$ws.Range("L94").Value = 1336.5385
$ws.Range("M94").Value = -33333583
$ws.Range("N94").Value = -2238.5385
# Row 99: H99=58824270, I99=90909496, K99=90909496, M99=-90907998
$ws.Range("H99").Value = 58824270
$ws.Range("I99").Value = 90909496
$ws.Range("K99").Value = 90909496
$ws.Range("M99").Value = -90907998

$ws = $wb.Worksheets.Item("CRP")
# Row 22: H22=724.1818, I22=103, K22=103, M22=247
$ws.Range("H22").Value = 724.1818
$ws.Range("I22").Value = 103
$ws.Range("K22").Value = 103
$ws.Range("M22").Value = 247
# Row 31: H31=32536.695, I31=970.5238000000001, J31=81639.63, K31=970.5238000000001, L31=81639.63, M31=-675.5238000000001, N31=-82229.63
$ws.Range("H31").Value = 32536.695
$ws.Range("I31").Value = 970.5238000000001
$ws.Range("J31").Value = 81639.63
$ws.Range("K31").Value = 970.5238000000001
$ws.Range("L31").Value = 81639.63
$ws.Range("M31").Value = -675.5238000000001
$ws.Range("N31").Value = -82229.63
# Row 34: H34=32536.695, I34=970.5238000000001, J34=81639.63, K34=970.5238000000001, L34=81639.63, M34=-768.5238000000001, N34=-82043.63
$ws.Range("H34").Value = 32536.695
$ws.Range("I34").Value = 970.5238000000001
$ws.Range("J34").Value = 81639.63
$ws.Range("K34").Value = 970.5238000000001
$ws.Range("L34").Value = 81639.63
$ws.Range("M34").Value = -768.5238000000001
$ws.Range("N34").Value = -82043.63
# Row 58: H58=989.14, I58=776.34283, J58=1485.6666, K58=776.34283, L58=1485.6666, M58=-573.34283, N58=-1891.6666
$ws.Range("H58").Value = 989.14
$ws.Range("I58").Value = 776.34283
$ws.Range("J58").Value = 1485.6666
$ws.Range("K58").Value = 776.34283
$ws.Range("L58").Value = 1485.6666
$ws.Range("M58").Value = -573.34283
$ws.Range("N58").Value = -1891.6666
# Row 132: H132=13336056, I132=18871242, J132=1288.2727, K132=56613726, L132=3864.8181, M132=-56611196, N132=-8924.8181
$ws.Range("H132").Value = 13336056
$ws.Range("I132").Value = 18871242
$ws.Range("J132").Value = 1288.2727
$ws.Range("K132").Value = 56613726
$ws.Range("L132").Value = 3864.8181
$ws.Range("M132").Value = -56611196
$ws.Range("N132").Value = -8924.8181
# Row 134: H134=2551.2585, I134=2769.02, J134=1190.25, K134=8307.059999999999, L134=3570.75, M134=-5772.059999999999, N134=-8640.75
$ws.Range("H134").Value = 2551.2585
$ws.Range("I134").Value = 2769.02
$ws.Range("J134").Value = 1190.25
$ws.Range("K134").Value = 8307.059999999999
$ws.Range("L134").Value = 3570.75
$ws.Range("M134").Value = -5772.059999999999
$ws.Range("N134").Value = -8640.75
# Row 136: H136=989.14, I136=776.34283, J136=1485.6666, K136=2329.02849, L136=4456.9998, M136=220.9715099999999, N136=-9556.9998
$ws.Range("H136").Value = 989.14
$ws.Range("I136").Value = 776.34283
$ws.Range("J136").Value = 1485.6666
$ws.Range("K136").Value = 2329.02849
$ws.Range("L136").Value = 4456.9998
$ws.Range("M136").Value = 220.9715099999999
$ws.Range("N136").Value = -9556.9998

$ws = $wb.Worksheets.Item("CUL")
# Row 5: H5=345570.8, I5=206.42857, J5=455459.5, K5=619.28571, L5=1366378.5, M5=-507.28571, N5=-1366602.5
$ws.Range("H5").Value = 345570.8
$ws.Range("I5").Value = 206.42857
$ws.Range("J5").Value = 455459.5
$ws.Range("K5").Value = 619.28571
$ws.Range("L5").Value = 1366378.5
$ws.Range("M5").Value = -507.28571
$ws.Range("N5").Value = -1366602.5
# Row 6: H6=172.78572, I6=34.75, J6=1001, K6=104.25, L6=3003, M6=8.75, N6=-3229
$ws.Range("H6").Value = 172.78572
$ws.Range("I6").Value = 34.75
$ws.Range("J6").Value = 1001
$ws.Range("K6").Value = 104.25
$ws.Range("L6").Value = 3003
$ws.Range("M6").Value = 8.75
$ws.Range("N6").Value = -3229
# Row 55: H55=2573.75, J55=3248.3333, L55=9744.999899999999, N55=-10098.9999
$ws.Range("H55").Value = 2573.75
$ws.Range("J55").Value = 3248.3333
$ws.Range("L55").Value = 9744.999899999999
$ws.Range("N55").Value = -10098.9999
# Row 131: H131=778.2033699999999, I131=454.95456, J131=970.4054, K131=1364.86368, L131=2911.2162, M131=3675.13632, N131=-12991.2162
$ws.Range("H131").Value = 778.2033699999999
$ws.Range("I131").Value = 454.95456
$ws.Range("J131").Value = 970.4054
$ws.Range("K131").Value = 1364.86368
$ws.Range("L131").Value = 2911.2162
$ws.Range("M131").Value = 3675.13632
$ws.Range("N131").Value = -12991.2162
# Row 134: H134=34369.83, I134=48186, J134=4225.4546, K134=144558, L134=12676.3638, M134=-139488, N134=-22816.3638
$ws.Range("H134").Value = 34369.83
$ws.Range("I134").Value = 48186
$ws.Range("J134").Value = 4225.4546
$ws.Range("K134").Value = 144558
$ws.Range("L134").Value = 12676.3638
$ws.Range("M134").Value = -139488
$ws.Range("N134").Value = -22816.3638
# Row 135: H135=345570.8, I135=206.42857, J135=455459.5, K135=1857.85713, L135=4099135.5, M135=677.1428699999999, N135=-4104205.5
$ws.Range("H135").Value = 345570.8
$ws.Range("I135").Value = 206.42857
$ws.Range("J135").Value = 455459.5
$ws.Range("K135").Value = 1857.85713
$ws.Range("L135").Value = 4099135.5
$ws.Range("M135").Value = 677.1428699999999
$ws.Range("N135").Value = -4104205.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22: H22=678.7059, I22=484.8, K22=484.8, M22=-189.8
$ws.Range("H22").Value = 678.7059
$ws.Range("I22").Value = 484.8
$ws.Range("K22").Value = 484.8
$ws.Range("M22").Value = -189.8
# Row 27: H27=678.7059, I27=484.8, K27=484.8, M27=-377.8
$ws.Range("H27").Value = 678.7059
$ws.Range("I27").Value = 484.8
$ws.Range("K27").Value = 484.8
$ws.Range("M27").Value = -377.8

$ws = $wb.Worksheets.Item("WVR")
# Row 130: H130=33694, J130=33694, L130=33694, N130=-43734
$ws.Range("H130").Value = 33694
$ws.Range("J130").Value = 33694
$ws.Range("L130").Value = 33694
$ws.Range("N130").Value = -43734
# Row 132: H132=2229.192, I132=2566.9814, J132=1269.1578, K132=7700.9442, L132=3807.4734, M132=-5170.9442, N132=-8867.473399999999
$ws.Range("H132").Value = 2229.192
$ws.Range("I132").Value = 2566.9814
$ws.Range("J132").Value = 1269.1578
$ws.Range("K132").Value = 7700.9442
$ws.Range("L132").Value = 3807.4734
$ws.Range("M132").Value = -5170.9442
$ws.Range("N132").Value = -8867.473399999999
